# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Updates the "K" column (column G) on the active worksheet with the
# recalculated values for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K (column G) value
$kValues = @{
    2  = 3
    3  = 3
    4  = 1
    7  = 2
    8  = 0
    9  = 1
    10 = 1
    11 = 0
    12 = 1
    13 = 0
    14 = 0
    15 = 1
    16 = 2
    17 = 1
    18 = 2
    20 = 0
    21 = 0
    22 = 1
    23 = 1
    24 = 3
    25 = 1
    26 = 0
    27 = 0
    28 = 1
    29 = 2
    30 = 2
    31 = 0
    32 = 0
    33 = 0
    34 = 1
    35 = 1
    36 = 0
    37 = 0
    38 = 0
    39 = 3
    40 = 1
    41 = 0
    42 = 1
    43 = 0
    44 = 2
    45 = 0
    46 = 0
    47 = 1
    48 = 1
    49 = 3
    50 = 0
    52 = 0
    54 = 1
    55 = 1
    56 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
